$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header count text
$ws.Range("B1").Value = "Count (Total: 235)"

# Update weekly triaged issue counts
$ws.Range("B2").Value = 136
$ws.Range("B3").Value = 83
$ws.Range("B4").Value = 16
